$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.533109716989713
$ws.Range("C2").Value = 0.1543646454752547
$ws.Range("D2").Value = 0.07743884050346139
$ws.Range("E2").Value = 0.07931802388768106
$ws.Range("G2").Value = 2.339728432799404
$ws.Range("H2").Value = 1.822285499500623
$ws.Range("L2").Value = 0.1929245253808247
$ws.Range("M2").Value = 0.3115730814082909
$ws.Range("B3").Value = 1.451471518257051
$ws.Range("C3").Value = 0.1378724350559821
$ws.Range("D3").Value = 0.07038185066704727
$ws.Range("E3").Value = 0.07924734858596338
$ws.Range("G3").Value = 2.270569127427819
$ws.Range("H3").Value = 1.795776888863912
$ws.Range("L3").Value = 0.1905517593952339
$ws.Range("M3").Value = 0.299297891348651
$ws.Range("B4").Value = 1.402432580623554
$ws.Range("C4").Value = 0.1276799121132228
$ws.Range("D4").Value = 0.06609154514812587
$ws.Range("E4").Value = 0.07920340269332382
$ws.Range("G4").Value = 2.229289719950401
$ws.Range("H4").Value = 1.780308884845141
$ws.Range("L4").Value = 0.1891946997608684
$ws.Range("M4").Value = 0.291974670912424
$ws.Range("B5").Value = 1.382721376698214
$ws.Range("C5").Value = 0.1235094406980295
$ws.Range("D5").Value = 0.06435377508535112
$ws.Range("E5").Value = 0.07918534789528953
$ws.Range("G5").Value = 2.212762921298946
$ws.Range("H5").Value = 1.774207552401094
$ws.Range("L5").Value = 0.1886667605471644
$ws.Range("M5").Value = 0.2890440336207192
$ws.Range("B6").Value = 1.379464784777781
$ws.Range("C6").Value = 0.1228159061945462
$ws.Range("D6").Value = 0.06406585262300268
$ws.Range("E6").Value = 0.07918234083124487
$ws.Range("G6").Value = 2.210036382411118
$ws.Range("H6").Value = 1.773206594378422
$ws.Range("L6").Value = 0.1885806102131298
$ws.Range("M6").Value = 0.2885606393703952
$ws.Range("B7").Value = 1.402165645930438
$ws.Range("C7").Value = 0.1276237366097916
$ws.Range("D7").Value = 0.0660680664281017
$ws.Range("E7").Value = 0.07920315980287718
$ws.Range("G7").Value = 2.229065643717945
$ws.Range("H7").Value = 1.780225783938789
$ws.Range("L7").Value = 0.1891874783024505
$ws.Range("M7").Value = 0.2919349303126921
$ws.Range("B8").Value = 1.504734663841703
$ws.Range("C8").Value = 0.1486917051171588
$ws.Range("D8").Value = 0.0749966158950599
$ws.Range("E8").Value = 0.07929376459233373
$ws.Range("G8").Value = 2.315634587036101
$ws.Range("H8").Value = 1.812976770723282
$ws.Range("L8").Value = 0.1920856571840233
$ws.Range("M8").Value = 0.3072961400253362
$ws.Range("B9").Value = 1.714551157517576
$ws.Range("C9").Value = 0.189494969528198
$ws.Range("D9").Value = 0.092853338048414
$ws.Range("E9").Value = 0.07946743185176131
$ws.Range("G9").Value = 2.494941826692639
$ws.Range("H9").Value = 1.883677835080732
$ws.Range("L9").Value = 0.1985629841359184
$ws.Range("M9").Value = 0.3391238050378576
$ws.Range("B10").Value = 1.874087377492174
$ws.Range("C10").Value = 0.2191848000686605
$ws.Range("D10").Value = 0.1061985205480909
$ws.Range("E10").Value = 0.07959309178150975
$ws.Range("G10").Value = 2.632715167387829
$ws.Range("H10").Value = 1.939663467977056
$ws.Range("L10").Value = 0.2038094517279916
$ws.Range("M10").Value = 0.3635610544735144
$ws.Range("B11").Value = 1.947855668173986
$ws.Range("C11").Value = 0.232634131775967
$ws.Range("D11").Value = 0.1123217188809065
$ws.Range("E11").Value = 0.0796499493889844
$ws.Range("G11").Value = 2.696750523913863
$ws.Range("H11").Value = 1.966031093588242
$ws.Range("L11").Value = 0.2063029207502751
$ws.Range("M11").Value = 0.3749103662924327
$ws.Range("B12").Value = 1.975963022450969
$ws.Range("C12").Value = 0.2377192657182263
$ws.Range("D12").Value = 0.114648176339216
$ws.Range("E12").Value = 0.07967144498094059
$ws.Range("G12").Value = 2.721198459571326
$ws.Range("H12").Value = 1.97614673323335
$ws.Range("L12").Value = 0.2072625489766438
$ws.Range("M12").Value = 0.3792417407652593
$ws.Range("B13").Value = 1.969901892172118
$ws.Range("C13").Value = 0.23662443446824
$ws.Range("D13").Value = 0.1141467848791109
$ws.Range("E13").Value = 0.07966681694263333
$ws.Range("G13").Value = 2.715924248992792
$ws.Range("H13").Value = 1.973962308993237
$ws.Range("L13").Value = 0.2070551899718538
$ws.Range("M13").Value = 0.378307402890897
$ws.Range("B14").Value = 1.950164605170414
$ws.Range("C14").Value = 0.2330526439896516
$ws.Range("D14").Value = 0.1125129617621354
$ws.Range("E14").Value = 0.07965171849777164
$ws.Range("G14").Value = 2.698757860049398
$ws.Range("H14").Value = 1.966860683732079
$ws.Range("L14").Value = 0.2063815609293584
$ws.Range("M14").Value = 0.3752660359142084
$ws.Range("B15").Value = 1.938097498833542
$ws.Range("C15").Value = 0.2308638083416952
$ws.Range("D15").Value = 0.1115132112261392
$ws.Range("E15").Value = 0.07964246595035429
$ws.Range("G15").Value = 2.68826898429819
$ws.Range("H15").Value = 1.962527809438257
$ws.Range("L15").Value = 0.2059709513062842
$ws.Range("M15").Value = 0.3734074966501026
$ws.Range("B16").Value = 1.869290578307016
$ws.Range("C16").Value = 0.2183047404424485
$ws.Range("D16").Value = 0.10579942665602
$ws.Range("E16").Value = 0.07958937060211424
$ws.Range("G16").Value = 2.628557989645401
$ws.Range("H16").Value = 1.937958504879759
$ws.Range("L16").Value = 0.2036486509970814
$ws.Range("M16").Value = 0.3628240473872424
$ws.Range("B17").Value = 1.827386488562126
$ws.Range("C17").Value = 0.2105858568837959
$ws.Range("D17").Value = 0.1023077684346561
$ws.Range("E17").Value = 0.07955672626029853
$ws.Range("G17").Value = 2.59227832881507
$ws.Range("H17").Value = 1.923117435960989
$ws.Range("L17").Value = 0.2022513862934829
$ws.Range("M17").Value = 0.3563911634029537
$ws.Range("B18").Value = 1.803396721833622
$ws.Range("C18").Value = 0.2061408007130296
$ws.Range("D18").Value = 0.1003043797765031
$ws.Range("E18").Value = 0.07953792078850719
$ws.Range("G18").Value = 2.571539183374739
$ws.Range("H18").Value = 1.914665821495106
$ws.Range("L18").Value = 0.2014577686606032
$ws.Range("M18").Value = 0.3527130290093936
$ws.Range("B19").Value = 1.795293462528718
$ws.Range("C19").Value = 0.2046348511131271
$ws.Range("D19").Value = 0.09962690594058188
$ws.Range("E19").Value = 0.07953154827841968
$ws.Range("G19").Value = 2.564539140262667
$ws.Range("H19").Value = 1.911818732165955
$ws.Range("L19").Value = 0.2011907884972572
$ws.Range("M19").Value = 0.35147143032966
$ws.Range("B20").Value = 1.831835613782744
$ws.Range("C20").Value = 0.2114080975056822
$ws.Range("D20").Value = 0.1026789506971824
$ws.Range("E20").Value = 0.07956020427607458
$ws.Range("G20").Value = 2.596127093579156
$ws.Range("H20").Value = 1.924688529370513
$ws.Range("L20").Value = 0.202399086840515
$ws.Range("M20").Value = 0.3570736880639203
$ws.Range("B21").Value = 1.955957222909547
$ws.Range("C21").Value = 0.2341019755105833
$ws.Range("D21").Value = 0.1129926436859279
$ws.Range("E21").Value = 0.07965615416028626
$ws.Range("G21").Value = 2.703794613305718
$ws.Range("H21").Value = 1.96894304193637
$ws.Range("L21").Value = 0.2065790036406696
$ws.Range("M21").Value = 0.3761584447906685
$ws.Range("B22").Value = 2.038086248755121
$ws.Range("C22").Value = 0.2488882163031292
$ws.Range("D22").Value = 0.1197784190673303
$ws.Range("E22").Value = 0.07971866224137081
$ws.Range("G22").Value = 2.775323852624922
$ws.Range("H22").Value = 1.998628964106331
$ws.Range("L22").Value = 0.2094006386239613
$ws.Range("M22").Value = 0.3888275827413423
$ws.Range("B23").Value = 1.994159785391673
$ws.Range("C23").Value = 0.2410005785294516
$ws.Range("D23").Value = 0.1161525224765541
$ws.Range("E23").Value = 0.07968531592806016
$ws.Range("G23").Value = 2.737039888581592
$ws.Range("H23").Value = 1.982714727875305
$ws.Range("L23").Value = 0.2078864454413889
$ws.Range("M23").Value = 0.3820478176125093
$ws.Range("B24").Value = 1.829823848353556
$ws.Range("C24").Value = 0.2110363854915249
$ws.Range("D24").Value = 0.1025111267139351
$ws.Range("E24").Value = 0.07955863198349578
$ws.Range("G24").Value = 2.594386698134429
$ws.Range("H24").Value = 1.923977986744234
$ws.Range("L24").Value = 0.2023322812490989
$ws.Range("M24").Value = 0.3567650557085997
$ws.Range("B25").Value = 1.656851641013589
$ws.Range("C25").Value = 0.1785089396930459
$ws.Range("D25").Value = 0.08798396574930223
$ws.Range("E25").Value = 0.0794208269759471
$ws.Range("G25").Value = 2.445390200901016
$ws.Range("H25").Value = 1.863848377369095
$ws.Range("L25").Value = 0.1967253336326849
$ws.Range("M25").Value = 0.3303296483190863
